# Refresh the cryptos list snapshot (prices + 1h volume deltas),
# matching the scheduled GitHub Actions data-refresh commit.
# Two pairs of rows also rotate: the coins that swapped ranking
# positions (rows 16/17, and the 3-way shuffle in rows 44-46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "60.546.60"
$ws.Range('E2').Value = "  +2.74%  "
$ws.Range('D3').Value = "2.705.82"
$ws.Range('E3').Value = "  +3.33%  "
$ws.Range('E4').Value = "  +0.04%  "
$ws.Range('D5').Value = "'526.89"
$ws.Range('E5').Value = "  +1.61%  "
$ws.Range('D6').Value = "'144.86"
$ws.Range('E6').Value = "  -0.22%  "
$ws.Range('E7').Value = "  +0.05%  "
$ws.Range('D8').Value = "'0.578"
$ws.Range('E8').Value = "  +2.13%  "
$ws.Range('D9').Value = "2.739.34"
$ws.Range('E9').Value = "  +3.83%  "
$ws.Range('D10').Value = "'6.70"
$ws.Range('E10').Value = "  +6.35%  "
$ws.Range('E11').Value = "  +1.07%  "
$ws.Range('E12').Value = "  +0.87%  "
$ws.Range('D14').Value = "3.187.05"
$ws.Range('E14').Value = "  +3.38%  "
$ws.Range('D15').Value = "60.575.88"
$ws.Range('E15').Value = "  +2.89%  "
$ws.Range('B16').Value = "WrappedEther"
$ws.Range('C16').Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('D16').Value = "2.848.79"
$ws.Range('E16').Value = "  +8.19%  "
$ws.Range('B17').Value = "Avalanche"
$ws.Range('C17').Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range('D17').Value = "'21.33"
$ws.Range('E17').Value = "  +1.71%  "
$ws.Range('E18').Value = "  +0.55%  "
$ws.Range('D19').Value = "'347.19"
$ws.Range('E19').Value = "  -0.57%  "
$ws.Range('D20').Value = "'4.52"
$ws.Range('E20').Value = "  +0.18%  "
$ws.Range('D21').Value = "'10.65"
$ws.Range('E21').Value = "  +3.31%  "
$ws.Range('E22').Value = "  +4.83%  "
$ws.Range('E23').Value = "  +0.10%  "
$ws.Range('D24').Value = "'63.71"
$ws.Range('E24').Value = "  +3.38%  "
$ws.Range('D25').Value = "'0.421"
$ws.Range('E25').Value = "  +0.58%  "
$ws.Range('E26').Value = "  +4.62%  "
$ws.Range('D27').Value = "'0.993"
$ws.Range('E27').Value = "  -0.07%  "
$ws.Range('D28').Value = "0.0₃0820"
$ws.Range('E28').Value = "  +2.06%  "
$ws.Range('E29').Value = "  +2.67%  "
$ws.Range('D30').Value = "'6.79"
$ws.Range('E30').Value = "  +9.05%  "
$ws.Range('D33').Value = "'19.14"
$ws.Range('E33').Value = "  +0.91%  "
$ws.Range('D34').Value = "'150.11"
$ws.Range('E34').Value = "  +0.12%  "
$ws.Range('E35').Value = "  +6.28%  "
$ws.Range('E36').Value = "  +8.27%  "
$ws.Range('E37').Value = "  -2.94%  "
$ws.Range('D38').Value = "'0.875"
$ws.Range('E38').Value = "  +3.88%  "
$ws.Range('E39').Value = "  +7.55%  "
$ws.Range('D40').Value = "'37.12"
$ws.Range('E40').Value = "  +1.23%  "
$ws.Range('D41').Value = "'3.67"
$ws.Range('E41').Value = "  -0.38%  "
$ws.Range('D42').Value = "'283.59"
$ws.Range('E42').Value = "  +2.27%  "
$ws.Range('D43').Value = "'20.19"
$ws.Range('E43').Value = "  +2.97%  "
$ws.Range('B44').Value = "Stellar"
$ws.Range('C44').Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range('D44').Value = "'0.0990"
$ws.Range('E44').Value = "  +0.56%  "
$ws.Range('B45').Value = "FirstDigitalUSD"
$ws.Range('C45').Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range('D45').Value = "'0.996"
$ws.Range('E45').Value = "  +0.17%  "
$ws.Range('B46').Value = "Mantle"
$ws.Range('C46').Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range('D46').Value = "'0.611"
$ws.Range('E46').Value = "  +0.90%  "
$ws.Range('D47').Value = "2.143.25"
$ws.Range('E47').Value = "  +8.28%  "
$ws.Range('D48').Value = "'0.0540"
$ws.Range('E48').Value = "  +3.27%  "
$ws.Range('D49').Value = "'4.83"
$ws.Range('E49').Value = "  +2.09%  "
$ws.Range('D50').Value = "'10.47"
$ws.Range('E50').Value = "  +1.80%  "
$ws.Range('E51').Value = "  +1.77%  "
